$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2228.8372
$ws.Range("I113").Value = 2030
$ws.Range("J113").Value = 2255
$ws.Range("K113").Value = 2030
$ws.Range("L113").Value = 2255
$ws.Range("M113").Value = 1224
$ws.Range("N113").Value = -8763

$ws.Range("H116").Value = 3866.6667
$ws.Range("I116").Value = 2460
$ws.Range("J116").Value = 4871.4287
$ws.Range("K116").Value = 2460
$ws.Range("L116").Value = 4871.4287
$ws.Range("M116").Value = 982
$ws.Range("N116").Value = -11755.4287

$ws.Range("H123").Value = 37413.5
$ws.Range("J123").Value = 37413.5
$ws.Range("L123").Value = 37413.5
$ws.Range("N123").Value = -47213.5

$ws.Range("H124").Value = 46659.25
$ws.Range("J124").Value = 46659.25
$ws.Range("L124").Value = 46659.25
$ws.Range("N124").Value = -56479.25

$ws.Range("H137").Value = 3272.1018
$ws.Range("I137").Value = 1043.025
$ws.Range("J137").Value = 7964.8945
$ws.Range("K137").Value = 3129.075
$ws.Range("L137").Value = 23894.6835
$ws.Range("M137").Value = -579.0750000000003
$ws.Range("N137").Value = -28994.6835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 45661.332
$ws.Range("J80").Value = 45661.332
$ws.Range("L80").Value = 45661.332
$ws.Range("N80").Value = -47657.332

$ws.Range("H83").Value = 45661.332
$ws.Range("J83").Value = 45661.332
$ws.Range("L83").Value = 136983.996
$ws.Range("N83").Value = -146967.996

$ws.Range("H128").Value = 50421
$ws.Range("J128").Value = 50421
$ws.Range("L128").Value = 50421
$ws.Range("N128").Value = -60381

$ws.Range("H139").Value = 46379.9
$ws.Range("J139").Value = 46379.9
$ws.Range("L139").Value = 46379.9
$ws.Range("N139").Value = -56659.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 37843.332
$ws.Range("J112").Value = 37843.332
$ws.Range("L112").Value = 37843.332
$ws.Range("N112").Value = -40797.332

$ws.Range("H122").Value = 41669
$ws.Range("J122").Value = 41669
$ws.Range("L122").Value = 41669
$ws.Range("N122").Value = -51469

$ws.Range("H132").Value = 40636
$ws.Range("J132").Value = 40636
$ws.Range("L132").Value = 40636
$ws.Range("N132").Value = -50756

$ws.Range("H133").Value = 47249.75
$ws.Range("J133").Value = 47249.75
$ws.Range("L133").Value = 47249.75
$ws.Range("N133").Value = -57369.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 43996
$ws.Range("J28").Value = 43996
$ws.Range("L28").Value = 43996
$ws.Range("N28").Value = -44486

$ws.Range("H43").Value = 48638.25
$ws.Range("J43").Value = 48638.25
$ws.Range("L43").Value = 48638.25
$ws.Range("N43").Value = -49006.25

$ws.Range("H95").Value = 110312
$ws.Range("J95").Value = 110312
$ws.Range("L95").Value = 110312
$ws.Range("N95").Value = -115804

$ws.Range("H96").Value = 42815.453
$ws.Range("J96").Value = 42815.453
$ws.Range("L96").Value = 42815.453
$ws.Range("N96").Value = -48307.453

$ws.Range("H100").Value = 40217.5
$ws.Range("J100").Value = 40217.5
$ws.Range("L100").Value = 40217.5
$ws.Range("N100").Value = -42381.5

$ws.Range("H101").Value = 48638.25
$ws.Range("J101").Value = 48638.25
$ws.Range("L101").Value = 48638.25
$ws.Range("N101").Value = -55128.25

$ws.Range("H106").Value = 48599.2
$ws.Range("J106").Value = 48599.2
$ws.Range("L106").Value = 48599.2
$ws.Range("N106").Value = -51123.2

$ws.Range("H112").Value = 39997.332
$ws.Range("J112").Value = 39997.332
$ws.Range("L112").Value = 39997.332
$ws.Range("N112").Value = -42951.332

$ws.Range("H118").Value = 48742
$ws.Range("J118").Value = 48742
$ws.Range("L118").Value = 48742
$ws.Range("N118").Value = -52056

$ws.Range("H134").Value = 312203.06
$ws.Range("I134").Value = 949.79486
$ws.Range("J134").Value = 2335349.2
$ws.Range("K134").Value = 2849.38458
$ws.Range("L134").Value = 7006047.600000001
$ws.Range("M134").Value = -314.3845799999999
$ws.Range("N134").Value = -7011117.600000001

$ws.Range("H137").Value = 42463.332
$ws.Range("J137").Value = 42463.332
$ws.Range("L137").Value = 42463.332
$ws.Range("N137").Value = -52663.332

$ws.Range("H139").Value = 74499.664
$ws.Range("J139").Value = 90749.5
$ws.Range("L139").Value = 90749.5
$ws.Range("N139").Value = -101029.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2273.125
$ws.Range("J34").Value = 2273.125
$ws.Range("L34").Value = 6819.375
$ws.Range("N34").Value = -6987.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 46665.25
$ws.Range("J104").Value = 46665.25
$ws.Range("L104").Value = 46665.25
$ws.Range("N104").Value = -53653.25

$ws.Range("H105").Value = 40790.8
$ws.Range("J105").Value = 40790.8
$ws.Range("L105").Value = 40790.8
$ws.Range("N105").Value = -47778.8

$ws.Range("H119").Value = 48761
$ws.Range("J119").Value = 48761
$ws.Range("L119").Value = 48761
$ws.Range("N119").Value = -58437

$ws.Range("H126").Value = 2616
$ws.Range("I126").Value = 3122.4
$ws.Range("J126").Value = 1350
$ws.Range("K126").Value = 9367.200000000001
$ws.Range("L126").Value = 4050
$ws.Range("M126").Value = -6897.200000000001
$ws.Range("N126").Value = -8990

$ws.Range("H133").Value = 48966.668
$ws.Range("J133").Value = 48966.668
$ws.Range("L133").Value = 48966.668
$ws.Range("N133").Value = -59086.668

$ws.Range("H135").Value = 37075.715
$ws.Range("J135").Value = 37075.715
$ws.Range("L135").Value = 37075.715
$ws.Range("N135").Value = -47215.715

$ws.Range("H138").Value = 51500
$ws.Range("J138").Value = 51500
$ws.Range("L138").Value = 51500
$ws.Range("N138").Value = -61780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 45387
$ws.Range("J111").Value = 45387
$ws.Range("L111").Value = 45387
$ws.Range("N111").Value = -53567

$ws.Range("H116").Value = 49676
$ws.Range("J116").Value = 49676
$ws.Range("L116").Value = 49676
$ws.Range("N116").Value = -58854

$ws.Range("H120").Value = 53129
$ws.Range("J120").Value = 53129
$ws.Range("L120").Value = 53129
$ws.Range("N120").Value = -62805

$ws.Range("H137").Value = 39966.668
$ws.Range("J137").Value = 39966.668
$ws.Range("L137").Value = 39966.668
$ws.Range("N137").Value = -50166.668

$ws.Range("H139").Value = 87666.336
$ws.Range("J139").Value = 51499.5
$ws.Range("L139").Value = 51499.5
$ws.Range("N139").Value = -61779.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45996
$ws.Range("J16").Value = 45996
$ws.Range("L16").Value = 45996
$ws.Range("N16").Value = -46580

$ws.Range("H110").Value = 48644
$ws.Range("J110").Value = 48644
$ws.Range("L110").Value = 48644
$ws.Range("N110").Value = -56824

$ws.Range("H138").Value = 42211.11
$ws.Range("J138").Value = 42211.11
$ws.Range("L138").Value = 42211.11
$ws.Range("N138").Value = -52491.11

$ws.Range("H139").Value = 58000
$ws.Range("J139").Value = 58000
$ws.Range("L139").Value = 58000
$ws.Range("N139").Value = -58854
